$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(448, "Pós OS", 8005278265, 46070.536539351851, "FRQ_ECO_PR_CURITIBA", "Detratores", "Ele é muito bom! Porém a mensalidade está se tornando um fardo! Esse ano troco por outro purificador!", "Outros", "Preço elevado"),
    @(449, "Pós OS", 8005285427, 46070.549305555563, "FRQ_ECO_RJ_OESTE", "Neutros", "Adoro o atendimento de vcs mais , desta vez nao trocou meu filtro.", "Campo", "Qualidade da manutenção"),
    @(450, "Pós OS", 8005283658, 46070.565648148149, "FRQ_ECO_SP_CAMPINAS_2", "Detratores", "O aparelho está fazendo barulho. Está atrapalhando dormir a noite.", "Qualidade do Produto", "Ruído"),
    @(451, "Pós OS", 8005285843, 46070.623692129629, "FRQ_ECO_SP_ZONA_SUL_03", "Neutros", "Normal", "Outros", "Satisfação geral"),
    @(452, "Pós OS", 8005293582, 46070.685486111113, "AT_ECO_SP_BAURU", "Neutros", "Água fresca e gostosa", "Outros", "Satisfação geral"),
    @(453, "Pós OS", 8005300015, 46071.391458333332, "FRQ_ECO_RJ_OESTE", "Detratores", "O pós-venda não está funcionando", "Outros", "Insatisfação geral"),
    @(454, "Pós OS", 8005304089, 46071.535034722219, "FRQ_ECO_SP_ZONA_SUL_03", "Detratores", "É a segunda vez que o técnico vem e não resolve o problema.Inicialmente trocou o filtro e depois de 3 dias o filtro apresentou o mesmo problema, daí o técnico voltou e disse que precisa substituir a peça.Substituir peça de um produto novo???? Estranho.E passa um prazo extenso para resolver, e preciso ligar para pedir prioridade.", "Campo", "Reincidência"),
    @(455, "Instalação", 8005295235, 46071.627256944441, "FRQ_ECO_RJ_OESTE", "Detratores", "Instalaram o equipamento, funcionou e no mesmo dia parou de funcionar a água gelada. Já solicitamos a presença do técnico há mais 03 dias e não fomos atendidos No aguardo", "Qualidade do Produto", "Não gela"),
    @(456, "Pós OS", 8005308413, 46072.391539351847, "FRQ_ECO_SP_CAMPINAS_2", "Detratores", "Mesmo solicitando um tecnico com urgencia, o purificador estava vazando agua, so agendaram 4 dias depois, disseram que foi por causa do carnaval, mas o vazamento não sai em bloquinho!!!!!!!", "Capacidade", "Agenda distante"),
    @(457, "Pós OS", 8005169116, 46072.397615740738, "FRQ_ECO_SP_ZONA_SUL_03", "Detratores", "Estamos a 3 semanas aguardando a troca do purificador que NÃO está gelando ..protocolo ( 8005280980 )..local que tem 30 pessoas tomando água quente 🥵", "Supply", "Falta de estoque/peças"),
    @(458, "Pós OS", 8005299637, 46072.496122685188, "FRQ_ECO_SP_ZONA_SUL_03", "Detratores", "Uma demora sem fim para o conserto", "Capacidade", "Data Distante"),
    @(459, "Instalação", 8005306854, 46072.659386574072, "FRQ_ECO_PR_CURITIBA", "Detratores", "Boa tarde muito mal instalado o puricador o rapaz nem ferramentas tinha para concluir o serviço solicito a reinstalação aqui ficou tudo feio a instalação", "Campo", "Qualidade da instalação"),
    @(460, "Instalação", 8005255862, 46073.40221064815, "FRQ_ECO_RJ_OESTE", "Detratores", "Péssima , até agora não funcionou, colocou um aparelho e só funcionou só um dia , já reclamei, marcaram uma visita hoje até às 13 horas , poderia dá uma preferência, obrigado.", "Qualidade do Produto", "Funcionamento geral"),
    @(461, "Pós OS", 8005304885, 46073.421550925923, "FRQ_ECO_SP_SANTOS", "Detratores", "Possuo a locação do purificador Brastemp a 10 anos e o produto é excelente .Infelizmente , pela primeira vez considerei em rescindir o contrato .Ocorre que estou com problemas no purificador sem poder usa-lo , já houveram duas visitas técnicas e , apesar do grande esforço dos técnicos , chegou-se à conclusão que seria necessário a troca do mesmo . A água gelada não funciona e nem sai ,  e a água natural sai muito quente . Apesar da situação grave já ocorrendo a trinta dias , sem poder utiliza-lo , a troca está marcada apenas para o dia 31/03/2026 . Pelo péssimo atendimento e negligência , o que senti começando a ocorrer recente , e apesar do produto ser excelente , ainda estou analisando a possibilidade de rescisão do contrato e compra de outro filtro de água . Inclusive , financeiramente seria vantajoso , pois  com o valor de 01 ano de locação conseguiria comprar um purificador de ótima qualidade .", "Supply", "Falta de estoque/peças"),
    @(462, "Pós OS", 8005304787, 46073.44326388889, "FRQ_ECO_RJ_ATLANTICA", "Neutros", "Boa. Duradouro", "Outros", "Satisfação geral"),
    @(463, "Pós OS", 8005293236, 46073.486342592587, "FRQ_ECO_SP_ZONA_SUL_03", "Detratores", "O purificador e muito bom e asistencia e boa mas o valor que e alto ainda mais que reformei o apartamebto e estou ha praticamente 2 anks morando fora e pagando sem usar o aparelho", "Outros", "Preço elevado"),
    @(464, "Pós OS", 8005286266, 46073.524375000001, "FRQ_ECO_RS_PORTO ALEGRE_2", "Detratores", "O purificador é bom, mas o sistema de marcação de visitas é deficiente", "Atendimento", "Erro de comunicação"),
    @(465, "Pós OS", 8005274003, 46073.552395833343, "FRQ_ECO_RJ_OESTE", "Detratores", "Muito frustrante o atendimento de manutenção. Minha instalação não foi concluída até hoje, apesar de várias remarcações para atendimento.", "Campo", "Qualidade da instalação"),
    @(466, "Pós OS", 8005286210, 46073.762187499997, "FRQ_ECO_SP_ZONA_SUL_03", "Neutros", "Apesar do conceito e do produto serem bons, a área de tecnologia, agendamento e solicitações de serviços é terrível. Não consigo imaginar algo pior. Agendamentos que não se completam, excesso de ligações de agendamento que falham, solicitações no agendamento ignoradas. Agendamentos duplicados. A área de TI é terrível", "Atendimento", "Qualidade do atendimento"),
    @(467, "Pós OS", 8005286288, 46073.766284722216, "FRQ_ECO_BA_SALVADOR", "Neutros", "Foi feito a troca dos filtros , mas ficou pingando o ponto de água na parede(mangueira) com a manipulação do aparelho. Deveria fazer um check list para tentar identificar esses problemas", "Campo", "Qualidade da manutenção"),
)

$startRow = 449
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

$ws.Range("G447").Select() | Out-Null
